$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date for the first row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-06 19:27:51"

# Sheet "zh-cn": Correspond Handoff Datetime / Correspond Handback DateTime (first row)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-06 19:27:46"
$wsZhCn.Range("K2").Value = "2016-09-06 19:28:12"

# Sheet "de-de": Correspond Handoff Datetime / Correspond Handback DateTime (first row)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-06 19:27:51"
$wsDeDe.Range("K2").Value = "2016-09-06 19:28:22"
